$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 19.484201272969916
$ws.Range("C2").Value = 5.4999095751223175
$ws.Range("D2").Value = 6.3435202707163434
$ws.Range("E2").Value = 1.783092086562087

$ws.Range("B3").Value = 33.143232035472344
$ws.Range("C3").Value = 4.2882762550519846
$ws.Range("D3").Value = -3.1778405313996139
$ws.Range("E3").Value = 9.4515205349522233

$ws.Range("B1:E3").Select()
